$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F2").Value = "YTR"
$ws.Range("G17").Select()
